$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # départements
$ws2 = $wb.Worksheets.Item(2)   # régions
$ws3 = $wb.Worksheets.Item(3)   # national

# --- départements ---
$ws1.Cells.Item(33, 3).Value = 31.95876288659793
$ws1.Cells.Item(33, 4).Value = 31
$ws1.Cells.Item(33, 5).Value = 97
$ws1.Cells.Item(53, 3).Value = 43.24324324324324
$ws1.Cells.Item(53, 4).Value = 16
$ws1.Cells.Item(59, 3).Value = 36.36363636363637
$ws1.Cells.Item(59, 4).Value = 24
$ws1.Cells.Item(59, 5).Value = 66
$ws1.Cells.Item(79, 3).Value = 43.95604395604396
$ws1.Cells.Item(79, 4).Value = 40
$ws1.Cells.Item(80, 3).Value = 51.35135135135135
$ws1.Cells.Item(80, 4).Value = 38
$ws1.Cells.Item(104, 3).Value = 3.90625
$ws1.Cells.Item(104, 4).Value = 5
$ws1.Cells.Item(121, 3).Value = 14.75409836065574
$ws1.Cells.Item(121, 4).Value = 9
$ws1.Cells.Item(130, 3).Value = 7.017543859649122
$ws1.Cells.Item(130, 4).Value = 12
$ws1.Cells.Item(132, 3).Value = 3.208556149732621
$ws1.Cells.Item(132, 4).Value = 6
$ws1.Cells.Item(146, 3).Value = 5.660377358490567
$ws1.Cells.Item(146, 4).Value = 3
$ws1.Cells.Item(146, 5).Value = 53
$ws1.Cells.Item(158, 3).Value = 7.109004739336493
$ws1.Cells.Item(158, 4).Value = 15
$ws1.Cells.Item(161, 3).Value = 6.25
$ws1.Cells.Item(161, 5).Value = 176
$ws1.Cells.Item(166, 3).Value = 6.4
$ws1.Cells.Item(166, 4).Value = 8
$ws1.Cells.Item(173, 3).Value = 12.22222222222222
$ws1.Cells.Item(173, 4).Value = 11
$ws1.Cells.Item(177, 3).Value = 7.407407407407407
$ws1.Cells.Item(177, 4).Value = 10
$ws1.Cells.Item(188, 3).Value = 10.3448275862069
$ws1.Cells.Item(188, 4).Value = 6
$ws1.Cells.Item(190, 3).Value = 6.299212598425196
$ws1.Cells.Item(190, 4).Value = 8
$ws1.Cells.Item(194, 3).Value = 5.769230769230769
$ws1.Cells.Item(194, 4).Value = 6
$ws1.Cells.Item(227, 3).Value = 5.232558139534884
$ws1.Cells.Item(227, 4).Value = 9
$ws1.Cells.Item(234, 3).Value = 4.878048780487805
$ws1.Cells.Item(234, 4).Value = 8
$ws1.Cells.Item(254, 3).Value = 6.122448979591836
$ws1.Cells.Item(254, 4).Value = 3
$ws1.Cells.Item(258, 3).Value = 4.205607476635514
$ws1.Cells.Item(258, 4).Value = 9
$ws1.Cells.Item(260, 3).Value = 3.96039603960396
$ws1.Cells.Item(260, 4).Value = 4
$ws1.Cells.Item(263, 3).Value = 5.839416058394161
$ws1.Cells.Item(263, 4).Value = 8
$ws1.Cells.Item(265, 3).Value = 7.602339181286549
$ws1.Cells.Item(265, 4).Value = 13
$ws1.Cells.Item(270, 3).Value = 6.666666666666667
$ws1.Cells.Item(270, 4).Value = 7
$ws1.Cells.Item(273, 3).Value = 3.529411764705882
$ws1.Cells.Item(273, 4).Value = 6
$ws1.Cells.Item(274, 3).Value = 3.546099290780142
$ws1.Cells.Item(274, 4).Value = 5
$ws1.Cells.Item(281, 3).Value = 7.07070707070707
$ws1.Cells.Item(281, 4).Value = 7
$ws1.Cells.Item(285, 3).Value = 12.82051282051282
$ws1.Cells.Item(285, 4).Value = 10
$ws1.Cells.Item(287, 3).Value = 5.263157894736842
$ws1.Cells.Item(287, 4).Value = 7
$ws1.Cells.Item(324, 3).Value = 2.538071065989848
$ws1.Cells.Item(324, 4).Value = 5
$ws1.Cells.Item(334, 3).Value = 1.612903225806452
$ws1.Cells.Item(334, 4).Value = 1
$ws1.Cells.Item(350, 3).Value = 3.90625
$ws1.Cells.Item(350, 4).Value = 5
$ws1.Cells.Item(359, 3).Value = 2.898550724637681
$ws1.Cells.Item(359, 4).Value = 2
$ws1.Cells.Item(369, 3).Value = 1.183431952662722
$ws1.Cells.Item(369, 4).Value = 2
$ws1.Cells.Item(370, 3).Value = 2.272727272727273
$ws1.Cells.Item(370, 4).Value = 4
$ws1.Cells.Item(371, 3).Value = 1.360544217687075
$ws1.Cells.Item(371, 4).Value = 2
$ws1.Cells.Item(381, 3).Value = 4.705882352941177
$ws1.Cells.Item(381, 4).Value = 4
$ws1.Cells.Item(388, 3).Value = 4.504504504504505
$ws1.Cells.Item(388, 4).Value = 5
$ws1.Cells.Item(395, 3).Value = 58.33333333333334
$ws1.Cells.Item(395, 4).Value = 21
$ws1.Cells.Item(525, 3).Value = 23.17073170731707
$ws1.Cells.Item(525, 4).Value = 19
$ws1.Cells.Item(534, 3).Value = 11.11111111111111
$ws1.Cells.Item(534, 4).Value = 4
$ws1.Cells.Item(534, 5).Value = 36
$ws1.Cells.Item(546, 3).Value = 17.42424242424243
$ws1.Cells.Item(546, 4).Value = 23
$ws1.Cells.Item(565, 3).Value = 23.52941176470588
$ws1.Cells.Item(565, 4).Value = 16
$ws1.Cells.Item(565, 5).Value = 68
$ws1.Cells.Item(615, 3).Value = 4.020100502512562
$ws1.Cells.Item(615, 4).Value = 8
$ws1.Cells.Item(651, 3).Value = 2.702702702702703
$ws1.Cells.Item(651, 4).Value = 4
$ws1.Cells.Item(653, 3).Value = 1.515151515151515
$ws1.Cells.Item(653, 4).Value = 3
$ws1.Cells.Item(661, 3).Value = 2.840909090909091
$ws1.Cells.Item(661, 4).Value = 5
$ws1.Cells.Item(662, 3).Value = 2.027027027027027
$ws1.Cells.Item(662, 4).Value = 3
$ws1.Cells.Item(673, 3).Value = 2.325581395348837
$ws1.Cells.Item(673, 4).Value = 2
$ws1.Cells.Item(712, 3).Value = 4.522613065326634
$ws1.Cells.Item(712, 4).Value = 9
$ws1.Cells.Item(722, 3).Value = 1.587301587301587
$ws1.Cells.Item(722, 4).Value = 1
$ws1.Cells.Item(738, 3).Value = 4.6875
$ws1.Cells.Item(738, 4).Value = 6
$ws1.Cells.Item(747, 3).Value = 2.898550724637681
$ws1.Cells.Item(747, 4).Value = 2
$ws1.Cells.Item(748, 3).Value = 4.026845637583892
$ws1.Cells.Item(748, 4).Value = 6
$ws1.Cells.Item(750, 3).Value = 1.515151515151515
$ws1.Cells.Item(750, 4).Value = 3
$ws1.Cells.Item(757, 3).Value = 1.764705882352941
$ws1.Cells.Item(757, 4).Value = 3
$ws1.Cells.Item(758, 3).Value = 2.824858757062147
$ws1.Cells.Item(758, 4).Value = 5
$ws1.Cells.Item(759, 3).Value = 2.702702702702703
$ws1.Cells.Item(759, 4).Value = 4
$ws1.Cells.Item(769, 3).Value = 4.705882352941177
$ws1.Cells.Item(769, 4).Value = 4
$ws1.Cells.Item(770, 3).Value = 2.325581395348837
$ws1.Cells.Item(770, 4).Value = 2
$ws1.Cells.Item(776, 3).Value = 6.194690265486726
$ws1.Cells.Item(776, 4).Value = 7
$ws1.Cells.Item(809, 3).Value = 2.010050251256281
$ws1.Cells.Item(809, 4).Value = 4
$ws1.Cells.Item(855, 3).Value = 2.259887005649718
$ws1.Cells.Item(855, 4).Value = 4
$ws1.Cells.Item(856, 3).Value = 0.6756756756756757
$ws1.Cells.Item(856, 4).Value = 1

# --- régions ---
$ws2.Cells.Item(2, 4).Value = 6.76
$ws2.Cells.Item(2, 5).Value = 55
$ws2.Cells.Item(3, 4).Value = 29.4
$ws2.Cells.Item(3, 5).Value = 122
$ws2.Cells.Item(3, 6).Value = 415
$ws2.Cells.Item(5, 4).Value = 2.4
$ws2.Cells.Item(5, 5).Value = 21
$ws2.Cells.Item(6, 4).Value = 55.34
$ws2.Cells.Item(6, 5).Value = 228
$ws2.Cells.Item(7, 4).Value = 5.84
$ws2.Cells.Item(7, 5).Value = 50
$ws2.Cells.Item(8, 4).Value = 3.97
$ws2.Cells.Item(8, 5).Value = 35
$ws2.Cells.Item(9, 4).Value = 1.81
$ws2.Cells.Item(9, 5).Value = 16
$ws2.Cells.Item(10, 4).Value = 3.41
$ws2.Cells.Item(10, 5).Value = 30
$ws2.Cells.Item(14, 4).Value = 0.45
$ws2.Cells.Item(14, 5).Value = 2
$ws2.Cells.Item(17, 4).Value = 0.67
$ws2.Cells.Item(17, 5).Value = 3
$ws2.Cells.Item(20, 4).Value = 9.699999999999999
$ws2.Cells.Item(20, 5).Value = 46
$ws2.Cells.Item(25, 4).Value = 8.91
$ws2.Cells.Item(25, 5).Value = 50
$ws2.Cells.Item(26, 4).Value = 2.26
$ws2.Cells.Item(26, 5).Value = 15
$ws2.Cells.Item(28, 4).Value = 2.13
$ws2.Cells.Item(28, 5).Value = 14
$ws2.Cells.Item(32, 4).Value = 1.6
$ws2.Cells.Item(32, 5).Value = 9
$ws2.Cells.Item(35, 4).Value = 1.95
$ws2.Cells.Item(35, 5).Value = 11
$ws2.Cells.Item(38, 4).Value = 6.55
$ws2.Cells.Item(38, 5).Value = 42
$ws2.Cells.Item(38, 6).Value = 641
$ws2.Cells.Item(39, 4).Value = 21.65
$ws2.Cells.Item(39, 5).Value = 84
$ws2.Cells.Item(43, 4).Value = 3.83
$ws2.Cells.Item(43, 5).Value = 28
$ws2.Cells.Item(47, 4).Value = 5.82
$ws2.Cells.Item(47, 5).Value = 41
$ws2.Cells.Item(50, 4).Value = 1.63
$ws2.Cells.Item(50, 5).Value = 14
$ws2.Cells.Item(51, 4).Value = 34.1
$ws2.Cells.Item(51, 5).Value = 148
$ws2.Cells.Item(51, 6).Value = 434
$ws2.Cells.Item(52, 4).Value = 4.78
$ws2.Cells.Item(52, 5).Value = 37
$ws2.Cells.Item(53, 4).Value = 2.21
$ws2.Cells.Item(53, 5).Value = 19
$ws2.Cells.Item(55, 4).Value = 1.63
$ws2.Cells.Item(55, 5).Value = 14
$ws2.Cells.Item(61, 4).Value = 4.84
$ws2.Cells.Item(61, 5).Value = 24
$ws2.Cells.Item(74, 4).Value = 6.23
$ws2.Cells.Item(74, 5).Value = 54
$ws2.Cells.Item(74, 6).Value = 867
$ws2.Cells.Item(75, 4).Value = 18.24
$ws2.Cells.Item(75, 5).Value = 79
$ws2.Cells.Item(75, 6).Value = 433
$ws2.Cells.Item(79, 4).Value = 3.67
$ws2.Cells.Item(79, 5).Value = 35
$ws2.Cells.Item(83, 4).Value = 7.33
$ws2.Cells.Item(83, 5).Value = 63
$ws2.Cells.Item(86, 4).Value = 1.38
$ws2.Cells.Item(86, 5).Value = 15
$ws2.Cells.Item(87, 4).Value = 26.48
$ws2.Cells.Item(87, 5).Value = 157
$ws2.Cells.Item(87, 6).Value = 593
$ws2.Cells.Item(88, 4).Value = 5.75
$ws2.Cells.Item(88, 5).Value = 52
$ws2.Cells.Item(89, 4).Value = 2.19
$ws2.Cells.Item(89, 5).Value = 24
$ws2.Cells.Item(90, 4).Value = 0.91
$ws2.Cells.Item(90, 5).Value = 10
$ws2.Cells.Item(91, 4).Value = 1.74
$ws2.Cells.Item(91, 5).Value = 19
$ws2.Cells.Item(92, 4).Value = 8.56
$ws2.Cells.Item(92, 5).Value = 87
$ws2.Cells.Item(93, 4).Value = 24.59
$ws2.Cells.Item(93, 5).Value = 135
$ws2.Cells.Item(97, 4).Value = 5.97
$ws2.Cells.Item(97, 5).Value = 66
$ws2.Cells.Item(98, 4).Value = 1.28
$ws2.Cells.Item(98, 5).Value = 17
$ws2.Cells.Item(100, 4).Value = 0.98
$ws2.Cells.Item(100, 5).Value = 13
$ws2.Cells.Item(101, 4).Value = 4.18
$ws2.Cells.Item(101, 5).Value = 29
$ws2.Cells.Item(103, 4).Value = 38.33
$ws2.Cells.Item(103, 5).Value = 69

# --- national ---
$ws3.Cells.Item(2, 2).Value = 6.97
$ws3.Cells.Item(2, 3).Value = 537
$ws3.Cells.Item(3, 2).Value = 22.19
$ws3.Cells.Item(3, 3).Value = 928
$ws3.Cells.Item(3, 4).Value = 4182
$ws3.Cells.Item(4, 2).Value = 46.21
$ws3.Cells.Item(4, 3).Value = 1190
$ws3.Cells.Item(5, 2).Value = 1.24
$ws3.Cells.Item(5, 3).Value = 121
$ws3.Cells.Item(6, 2).Value = 28.71
$ws3.Cells.Item(6, 3).Value = 1524
$ws3.Cells.Item(6, 4).Value = 5309
$ws3.Cells.Item(7, 2).Value = 5.2
$ws3.Cells.Item(7, 3).Value = 442
$ws3.Cells.Item(8, 2).Value = 1.73
$ws3.Cells.Item(8, 3).Value = 170
$ws3.Cells.Item(9, 2).Value = 0.91
$ws3.Cells.Item(9, 3).Value = 89
$ws3.Cells.Item(10, 2).Value = 1.41
$ws3.Cells.Item(10, 3).Value = 138
